$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C2").Value = 0.2314448356628418
$ws.Range("E2").Value = 135.377112823955
$ws.Range("F2").Value = 0.004501256144425741
$ws.Range("G2").Value = 0.003975177641773857
$ws.Range("H2").Value = 0.003794857390920524
$ws.Range("I2").Value = 0.003374893805232832
$ws.Range("J2").Value = 0.003153323868414926
$ws.Range("K2").Value = 0.003153323868414926
$ws.Range("L2").Value = 0.003053946330232168
$ws.Range("M2").Value = 0.003053946330232168
$ws.Range("N2").Value = 0.003053946330232168
$ws.Range("O2").Value = 0.002915008256453515
$ws.Range("P2").Value = 0.002915008256453515
$ws.Range("Q2").Value = 0.002915008256453515
$ws.Range("R2").Value = 0.002847843814783699
$ws.Range("S2").Value = 0.002738634971041211
$ws.Range("T2").Value = 0.002713189473323647
$ws.Range("U2").Value = 0.002713189473323647
$ws.Range("V2").Value = 0.002692047660081536
$ws.Range("W2").Value = 0.00264321649884499
$ws.Range("X2").Value = 0.002642589193198934
$ws.Range("Y2").Value = 0.002638930074541033

$ws.Range("C3").Value = 0.3429763317108154
$ws.Range("E3").Value = 140.6402111515235
$ws.Range("F3").Value = 0.004615634665128008
$ws.Range("G3").Value = 0.003950969851993556
$ws.Range("H3").Value = 0.003488287948609589
$ws.Range("I3").Value = 0.003488287948609589
$ws.Range("J3").Value = 0.003335573532809154
$ws.Range("K3").Value = 0.003335573532809154
$ws.Range("L3").Value = 0.003254173971923506
$ws.Range("M3").Value = 0.003173396093145029
$ws.Range("N3").Value = 0.003173396093145029
$ws.Range("O3").Value = 0.003110260335923912
$ws.Range("P3").Value = 0.003110260335923912
$ws.Range("Q3").Value = 0.002938858922996793
$ws.Range("R3").Value = 0.002938858922996793
$ws.Range("S3").Value = 0.002938858922996793
$ws.Range("T3").Value = 0.002910469537963183
$ws.Range("U3").Value = 0.002890734417561941
$ws.Range("V3").Value = 0.002878008044052987
$ws.Range("W3").Value = 0.00279214166407306
$ws.Range("X3").Value = 0.00279214166407306
$ws.Range("Y3").Value = 0.002741524583850361

$ws.Range("C4").Value = 0.3399138450622559
$ws.Range("E4").Value = 137.3039868001506
$ws.Range("F4").Value = 0.00433510701669386
$ws.Range("G4").Value = 0.003736444096399674
$ws.Range("H4").Value = 0.003658445665839916
$ws.Range("I4").Value = 0.003425334252635089
$ws.Range("J4").Value = 0.003238710785821449
$ws.Range("K4").Value = 0.003238710785821449
$ws.Range("L4").Value = 0.003238710785821449
$ws.Range("M4").Value = 0.003238710785821449
$ws.Range("N4").Value = 0.003207013671785601
$ws.Range("O4").Value = 0.00298866002835445
$ws.Range("P4").Value = 0.00298866002835445
$ws.Range("Q4").Value = 0.00298866002835445
$ws.Range("R4").Value = 0.00294705848722366
$ws.Range("S4").Value = 0.002922935500142631
$ws.Range("T4").Value = 0.002905910663629571
$ws.Range("U4").Value = 0.002809266476760967
$ws.Range("V4").Value = 0.002769861894226586
$ws.Range("W4").Value = 0.002725667021769805
$ws.Range("X4").Value = 0.002701392159600034
$ws.Range("Y4").Value = 0.002676490970763169

$ws.Range("C5").Value = 0.2355992794036865
$ws.Range("E5").Value = 138.5499295047157
$ws.Range("F5").Value = 0.004398872970596697
$ws.Range("G5").Value = 0.003835961736963654
$ws.Range("H5").Value = 0.003523989695630142
$ws.Range("I5").Value = 0.003523989695630142
$ws.Range("J5").Value = 0.003511824940256811
$ws.Range("K5").Value = 0.003469151577733073
$ws.Range("L5").Value = 0.003220804998562735
$ws.Range("M5").Value = 0.003179109999931803
$ws.Range("N5").Value = 0.003179109999931803
$ws.Range("O5").Value = 0.003030883081295184
$ws.Range("P5").Value = 0.002938533738268917
$ws.Range("Q5").Value = 0.002938533738268917
$ws.Range("R5").Value = 0.002835127482356514
$ws.Range("S5").Value = 0.002835127482356514
$ws.Range("T5").Value = 0.002775367609146837
$ws.Range("U5").Value = 0.002775367609146837
$ws.Range("V5").Value = 0.002775367609146837
$ws.Range("W5").Value = 0.002753907256612779
$ws.Range("X5").Value = 0.002741843535774842
$ws.Range("Y5").Value = 0.002700778352918434

$ws.Range("C6").Value = 0.3371186256408691
$ws.Range("E6").Value = 137.7672016294746
$ws.Range("F6").Value = 0.004627253007860288
$ws.Range("G6").Value = 0.003718800522516498
$ws.Range("H6").Value = 0.003718800522516498
$ws.Range("I6").Value = 0.003717636637804077
$ws.Range("J6").Value = 0.003500037738349082
$ws.Range("K6").Value = 0.003500037738349082
$ws.Range("L6").Value = 0.003349339282146305
$ws.Range("M6").Value = 0.003343842457460069
$ws.Range("N6").Value = 0.003143097471198863
$ws.Range("O6").Value = 0.003143097471198863
$ws.Range("P6").Value = 0.003098251838233855
$ws.Range("Q6").Value = 0.003053454012777286
$ws.Range("R6").Value = 0.003001205619705244
$ws.Range("S6").Value = 0.002947765901980883
$ws.Range("T6").Value = 0.002890653231714214
$ws.Range("U6").Value = 0.002880827882316785
$ws.Range("V6").Value = 0.00276398684989353
$ws.Range("W6").Value = 0.002726578528393671
$ws.Range("X6").Value = 0.002705271352275289
$ws.Range("Y6").Value = 0.002685520499599894

$ws.Range("C7").Value = 0.3339564800262451
$ws.Range("E7").Value = 138.3781393995232
$ws.Range("F7").Value = 0.004447304750687902
$ws.Range("G7").Value = 0.003956789097079046
$ws.Range("H7").Value = 0.003731461609205266
$ws.Range("I7").Value = 0.003527698161883749
$ws.Range("J7").Value = 0.003260446624240572
$ws.Range("K7").Value = 0.003260446624240572
$ws.Range("L7").Value = 0.003260446624240572
$ws.Range("M7").Value = 0.003188012404131447
$ws.Range("N7").Value = 0.003048376225687655
$ws.Range("O7").Value = 0.00302014738982998
$ws.Range("P7").Value = 0.002966249496182418
$ws.Range("Q7").Value = 0.002966249496182418
$ws.Range("R7").Value = 0.002966249496182418
$ws.Range("S7").Value = 0.002831872007076026
$ws.Range("T7").Value = 0.002831872007076026
$ws.Range("U7").Value = 0.002828657232220806
$ws.Range("V7").Value = 0.002757586714101529
$ws.Range("W7").Value = 0.002757586714101529
$ws.Range("X7").Value = 0.002721418460848023
$ws.Range("Y7").Value = 0.002697429617924429

$ws.Range("C8").Value = 0.2606539726257324
$ws.Range("E8").Value = 139.3758696982768
$ws.Range("F8").Value = 0.00447304128724669
$ws.Range("G8").Value = 0.003961682878800228
$ws.Range("H8").Value = 0.003648478404328215
$ws.Range("I8").Value = 0.003643613416052122
$ws.Range("J8").Value = 0.003571207560830744
$ws.Range("K8").Value = 0.003571041832015219
$ws.Range("L8").Value = 0.003494837615210689
$ws.Range("M8").Value = 0.003275609430618193
$ws.Range("N8").Value = 0.003184412653084075
$ws.Range("O8").Value = 0.003099132162155157
$ws.Range("P8").Value = 0.003099132162155157
$ws.Range("Q8").Value = 0.003049876964437272
$ws.Range("R8").Value = 0.002847593818675844
$ws.Range("S8").Value = 0.002847593818675844
$ws.Range("T8").Value = 0.002846794605923075
$ws.Range("U8").Value = 0.002808538268293498
$ws.Range("V8").Value = 0.002780513404557116
$ws.Range("W8").Value = 0.002780513404557116
$ws.Range("X8").Value = 0.002738476924051853
$ws.Range("Y8").Value = 0.002716878551623329

$ws.Range("C9").Value = 0.2480313777923584
$ws.Range("E9").Value = 136.6947743276414
$ws.Range("F9").Value = 0.004380533545867449
$ws.Range("G9").Value = 0.003826871150732757
$ws.Range("H9").Value = 0.003595881178480978
$ws.Range("I9").Value = 0.003595881178480978
$ws.Range("J9").Value = 0.003595881178480978
$ws.Range("K9").Value = 0.003451007987261537
$ws.Range("L9").Value = 0.003399213325043858
$ws.Range("M9").Value = 0.00324767354144711
$ws.Range("N9").Value = 0.00316906211958885
$ws.Range("O9").Value = 0.003168729533410324
$ws.Range("P9").Value = 0.002941680590475941
$ws.Range("Q9").Value = 0.002941680590475941
$ws.Range("R9").Value = 0.002882847382758224
$ws.Range("S9").Value = 0.002870971220770699
$ws.Range("T9").Value = 0.002811638296305133
$ws.Range("U9").Value = 0.002796849816869958
$ws.Range("V9").Value = 0.002749896561195485
$ws.Range("W9").Value = 0.002686855478302044
$ws.Range("X9").Value = 0.002686855478302044
$ws.Range("Y9").Value = 0.002664615483969618

$ws.Range("C10").Value = 0.2206225395202637
$ws.Range("E10").Value = 138.3581620700224
$ws.Range("F10").Value = 0.004465899194877027
$ws.Range("G10").Value = 0.003907829559908991
$ws.Range("H10").Value = 0.003511412453748877
$ws.Range("I10").Value = 0.003339435791055122
$ws.Range("J10").Value = 0.003162764901593703
$ws.Range("K10").Value = 0.00315728069415709
$ws.Range("L10").Value = 0.00315728069415709
$ws.Range("M10").Value = 0.00315728069415709
$ws.Range("N10").Value = 0.00315728069415709
$ws.Range("O10").Value = 0.003044564663160718
$ws.Range("P10").Value = 0.002987461824108336
$ws.Range("Q10").Value = 0.002878739011244427
$ws.Range("R10").Value = 0.002806607257657179
$ws.Range("S10").Value = 0.002806607257657179
$ws.Range("T10").Value = 0.002785732892269853
$ws.Range("U10").Value = 0.00278482731724009
$ws.Range("V10").Value = 0.002763556253634319
$ws.Range("W10").Value = 0.002763556253634319
$ws.Range("X10").Value = 0.002721242595486088
$ws.Range("Y10").Value = 0.002697040196296733

$ws.Range("C11").Value = 0.2084298133850098
$ws.Range("E11").Value = 136.7230257424308
$ws.Range("F11").Value = 0.004608920577053108
$ws.Range("G11").Value = 0.003607200726702097
$ws.Range("H11").Value = 0.003607200726702097
$ws.Range("I11").Value = 0.003607200726702097
$ws.Range("J11").Value = 0.003394792618068765
$ws.Range("K11").Value = 0.003394792618068765
$ws.Range("L11").Value = 0.003385472471467504
$ws.Range("M11").Value = 0.003353271968289225
$ws.Range("N11").Value = 0.003123160999319678
$ws.Range("O11").Value = 0.002976491812094831
$ws.Range("P11").Value = 0.002852683388897261
$ws.Range("Q11").Value = 0.002852683388897261
$ws.Range("R11").Value = 0.002852683388897261
$ws.Range("S11").Value = 0.002802443494147793
$ws.Range("T11").Value = 0.002797224871341398
$ws.Range("U11").Value = 0.002743100358536828
$ws.Range("V11").Value = 0.002743100358536828
$ws.Range("W11").Value = 0.002682333297161711
$ws.Range("X11").Value = 0.002682333297161711
$ws.Range("Y11").Value = 0.002665166193809566
